$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "삶의질지수"
$ws.Range("D1").Value = "양호한주관적건강수준인지율"
$ws.Range("E1").Value = "스트레스인지율"
$ws.Range("F1").Value = "우울감경험률"
$ws.Range("G1").Value = "주관적구강건강이나쁜인구의분율"
$ws.Range("I1").Value = "연간보건기관이용률"
